# Adds a new "Rules Visualization with arulesViz" slide at the end of the
# deck (slide 4), using the same "Title and Content" layout as the other
# slides in the presentation.

$p = $ppt.ActivePresentation

# ppLayoutText (2) maps to the "Title and Content" custom layout already
# used by every other slide in this deck (a title placeholder + a single
# body/content placeholder).
$s = $p.Slides.Add($p.Slides.Count + 1, 2)

# Title: "Rules Visualization with arulesViz" (kept as two runs so the
# product name stands apart from the rest of the sentence).
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Rules Visualization with "
[void]$title.InsertAfter("arulesViz")

# Body / content placeholder bullets.
$body = $s.Shapes.Item(2).TextFrame.TextRange
$lines = @(
    "37 rules : min support = 1%, min confidence = 30%",
    "Visualization types:",
    "Scatter plot",
    "Matrix (3DMatrix)",
    "Grouped",
    "Graph",
    "Parallel coordinate"
)
$body.Text = [string]::Join("`r", $lines)

# Second-level (indented) bullets for the visualization type list.
for ($i = 3; $i -le 7; $i++) {
    $body.Paragraphs($i).IndentLevel = 2
}
